$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.378.84'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.721.65'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''244.10'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +1.97%  '
$ws.Range("D8").Value = '''0.2608'
$ws.Range("E8").Value = '  -2.25%  '
$ws.Range("D9").Value = '''0.06204'
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("D10").Value = '1.727.06'
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("D11").Value = '''0.07013'
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("D12").Value = '''15.41'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '''4.531'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").Value = '''0.5967'
$ws.Range("E14").Value = '  -2.85%  '
$ws.Range("D15").Value = '''77.11'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '26.389.22'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").Value = '''1.001'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '''0.000007200'
$ws.Range("E19").Value = '  +3.19%  '
$ws.Range("D20").Value = '''11.35'
$ws.Range("E20").Value = '  -2.46%  '
$ws.Range("D21").Value = '1.949.92'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = '''4.474'
$ws.Range("E22").Value = '  -1.12%  '
$ws.Range("D23").Value = '''8.549'
$ws.Range("E23").Value = '  -4.24%  '
$ws.Range("D24").Value = '''5.167'
$ws.Range("E24").Value = '  -2.23%  '
$ws.Range("D25").Value = '''137.35'
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("D26").Value = '''15.21'
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").Value = '''1.410'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").Value = '''107.29'
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").Value = '''1.712'
$ws.Range("E29").Value = '  -4.44%  '
$ws.Range("D30").Value = '''3.953'
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").Value = '''0.07949'
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").Value = '''3.671'
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("D33").Value = '''0.04518'
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").Value = '''1.000'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("D36").Value = '''0.9927'
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").Value = '''0.6205'
$ws.Range("E37").Value = '  -1.71%  '
$ws.Range("D38").Value = '''0.9060'
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("D39").Value = '''1.977'
$ws.Range("E39").Value = '  -5.60%  '
$ws.Range("D40").Value = '''2.395'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("D42").Value = '''0.01485'
$ws.Range("E42").Value = '  -1.22%  '
$ws.Range("D43").Value = '''100.29'
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("D44").Value = '''5.394'
$ws.Range("E44").Value = '  -3.39%  '
$ws.Range("D45").Value = '''0.3840'
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("D46").Value = '''6.720'
$ws.Range("E46").Value = '  -3.47%  '
$ws.Range("D47").Value = '''0.1147'
$ws.Range("E47").Value = '  -3.16%  '
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''30.11'
$ws.Range("E49").Value = '  -2.68%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''7.675'
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").Value = '''1.242'
$ws.Range("E51").Value = '  -1.29%  '
